# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# immediately before the existing "Late" column (column N), pushing the
# "Late", "Date"(heading) and "Outstanding" columns one slot to the
# right (N->O, O->P, P->Q). The new column inherits its width from the
# column immediately to its left ("In Advance", column M), which is
# what Excel does when a column is inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late").
$ws.Range("N1").EntireColumn.Insert()

# The newly inserted column picks up the column width of its left
# neighbour (column M, "In Advance"), matching Excel's default
# insert-column behaviour.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Update the active cell/selection on the sheet to match the edited
# workbook.
$ws.Range("J20").Select() | Out-Null
